$d = $word.ActiveDocument

function Replace-Text([string]$findText, [string]$replaceText) {
    $r = $d.Content
    $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

function Find-Text([string]$findText) {
    $r = $d.Content
    $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r
}

function Get-ParagraphContaining($pos) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $pos -and $pos -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

function Get-ParagraphForText([string]$text) {
    $r = Find-Text $text
    return Get-ParagraphContaining($r.Start)
}

# --- 1. Plain text fixups (these merge split runs and/or change wording) ---

# "M" + "ust have permission to install plug ins." -> single run, same text
Replace-Text "Must have permission to install plug ins." "Must have permission to install plug ins."

# "M" + "onitor" + "s" + " user account behavior through the administrative dashboard" -> single run
# (leave the trailing "." in its own run, untouched)
Replace-Text "Monitors user account behavior through the administrative dashboard" "Monitors user account behavior through the administrative dashboard"

# "L" + "ogin through the administrative login page" + "." -> single run (period merges this time)
Replace-Text "Login through the administrative login page." "Login through the administrative login page."

# " ...and hit the landing " + "page." -> single run
Replace-Text " to view the website via domain name on their web browser of choice and hit the landing page." " to view the website via domain name on their web browser of choice and hit the landing page."

# " to visit permitted pages on the " + "website." -> single run
Replace-Text " to visit permitted pages on the website." " to visit permitted pages on the website."

# "Users cannot edit any pages on the site." -> "User cannot edit any pages on the site."
Replace-Text "Users cannot edit any pages on the site." "User cannot edit any pages on the site."

# "Users can fill out the contact us page." -> "User can fill out the contact us page."
Replace-Text "Users can fill out the contact us page." "User can fill out the contact us page."

# "Keeps all users on same page, no redirects." -> "...no redirects except for payments."
Replace-Text "Keeps all users on same page, no redirects." "Keeps all users on same page, no redirects except for payments."

# --- 2. Remove the "Testimonial page allowing for reviews." bullet entirely ---

$r = Find-Text "Testimonial page allowing for reviews."
if ($r.Find.Found) {
    $r.Expand(4)
    $r.Delete()
}

# --- 3. Add two new bullets under "Unauthenticated User" list (numId=1), after the
#        "User can fill out the contact us page." item ---

$p = Get-ParagraphForText "User can fill out the contact us page."
$p.Range.InsertParagraphAfter()

$p = Get-ParagraphForText "User can fill out the contact us page."
$newPara = $p.Next()
$newPara.Range.InsertBefore("User can fill out a testimonial.")
$newPara.Range.InsertParagraphAfter()

$p = Get-ParagraphForText "User can fill out a testimonial."
$newPara2 = $p.Next()
$newPara2.Range.InsertBefore("User has ability to make payments for services.")

# --- 4. Add two new bullets under the Functional Requirements list (numId=4), after
#        "Security check/ reCAPTCHA" ---

$p = Get-ParagraphForText "Security check/ reCAPTCHA"
$p.Range.InsertParagraphAfter()

$p = Get-ParagraphForText "Security check/ reCAPTCHA"
$newPara3 = $p.Next()
$newPara3.Range.InsertBefore("Payment page.")
$newPara3.Range.InsertParagraphAfter()

$p = Get-ParagraphForText "Payment page."
$newPara4 = $p.Next()
$newPara4.Range.InsertBefore("Testimonial page.")
